# Apply the "Atualização de bases das ligas" edit:
# 1) Several match rows had been listed with swapped content relative to
#    their correct chronological/ID order; restore each pair's data.
# 2) Append the new fixtures (rows 96-102) that were added to the bottom
#    of the sheet, copying formatting from the last existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-pair corrections (rows whose match data were in the wrong order) ---
$ws.Range("B4").Value2 = 6832698
$row = @("Al Ain SCC", "Baniyas SC", 3, 2, "H", 1.25, 6.5, 7.5, 1.25, 6.5, 7.5, -1.75, 1.825, 1.975, 3.5, 2, 1.8, 0.25, -1, -1, -1, 0.9750000000000001, 1, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F4:AC4").Value2 = $arr

$ws.Range("B5").Value2 = 6832494
$row = @("Ajman SCC", "Shabab Al Ahli Dubai", 0, 3, "A", 4.2, 3.6, 1.727, 5.25, 4, 1.533, 1, 1.925, 1.875, 3, 1.9, 1.9, -1, -1, 0.5329999999999999, -1, 0.875, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F5:AC5").Value2 = $arr

$ws.Range("B11").Value2 = 6832497
$row = @("Al Bataeh", "Al Ittihad Kalba", 0, 0, "D", 3.25, 3.3, 2.05, 3.5, 3.4, 1.95, 0.5, 1.8, 2, 2.75, 1.8, 2, -1, 2.4, -1, 0.8, -1, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F11:AC11").Value2 = $arr

$ws.Range("B12").Value2 = 6832699
$row = @("Hatta Dubai", "Al Ain SCC", 0, 2, "A", 5, 4.333, 1.5, 7, 5, 1.3, 1.5, 1.975, 1.825, 3.25, 1.95, 1.85, -1, -1, 0.3, -1, 0.825, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F12:AC12").Value2 = $arr

$ws.Range("B16").Value2 = 6832504
$row = @("Al Ittihad Kalba", "Al Jazira SC", 0, 4, "A", 3.6, 3.8, 1.85, 3.8, 4, 1.727, 0.75, 1.875, 1.925, 3.25, 1.95, 1.85, -1, -1, 0.7270000000000001, -1, 0.925, 0.95, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F16:AC16").Value2 = $arr

$ws.Range("B17").Value2 = 6832502
$row = @("Al Nasr SC", "Sharjah SCC", 0, 1, "A", 3.75, 3.6, 1.909, 3.8, 3.8, 1.8, 0.5, 2, 1.8, 3, 2, 1.8, -1, -1, 0.8, -1, 0.8, -1, 0.8)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F17:AC17").Value2 = $arr

$ws.Range("B25").Value2 = 6832512
$row = @("Ajman SCC", "Al Ittihad Kalba", 3, 5, "A", 2.4, 3.5, 2.55, 2.875, 3.6, 2.15, 0.25, 1.85, 1.95, 3, 1.9, 1.9, -1, -1, 1.15, -1, 0.95, 0.8999999999999999, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F25:AC25").Value2 = $arr

$ws.Range("B26").Value2 = 6832701
$row = @("Baniyas SC", "Emirates Club RAK", 0, 0, "D", 1.666, 4, 4, 1.8, 4, 3.5, -0.5, 1.825, 1.975, 3.25, 1.925, 1.875, -1, 3, -1, -1, 0.9750000000000001, -1, 0.875)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F26:AC26").Value2 = $arr

$ws.Range("B63").Value2 = 6832544
$row = @("Al Jazira SC", "Ajman SCC", 1, 5, "A", 1.222, 6.5, 11, 1.4, 5, 6.5, -1.5, 1.95, 1.85, 3.75, 1.95, 1.85, -1, -1, 5.5, -1, 0.8500000000000001, 0.95, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F63:AC63").Value2 = $arr

$ws.Range("B64").Value2 = 6832705
$row = @("Al Bataeh", "Baniyas SC", 2, 1, "H", 2.6, 3.4, 2.6, 2.25, 3.4, 3.1, -0.25, 2, 1.8, 2.75, 1.85, 1.95, 1.25, -1, -1, 1, -1, 0.425, -0.5)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F64:AC64").Value2 = $arr

$ws.Range("B67").Value2 = 6832549
$row = @("Hatta Dubai", "Sharjah SCC", 1, 4, "A", 15, 7, 1.111, 10, 5, 1.25, 1.75, 1.85, 1.95, 3, 1.85, 1.95, -1, -1, 0.25, -1, 0.95, 0.8500000000000001, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F67:AC67").Value2 = $arr

$ws.Range("B68").Value2 = 6832547
$row = @("Khor Fakkan", "Al Nasr SC", 3, 1, "H", 2.75, 3.4, 2.45, 6, 4.333, 1.5, 1, 2, 1.8, 3, 1.925, 1.875, 5, -1, -1, 1, -1, 0.925, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F68:AC68").Value2 = $arr

$ws.Range("B70").Value2 = 6832706
$row = @("Hatta Dubai", "Al Bataeh", 0, 0, "D", 3.3, 3.6, 1.909, 4, 3.75, 1.727, 0.75, 1.825, 1.975, 3, 2, 1.8, -1, 2.75, -1, 0.825, -1, -1, 0.8)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F70:AC70").Value2 = $arr

$ws.Range("B71").Value2 = 6832707
$row = @("Emirates Club RAK", "Al Jazira SC", 0, 1, "A", 5.75, 4.75, 1.4, 4.75, 4.5, 1.533, 1, 2, 1.8, 4, 2, 1.8, -1, -1, 0.5329999999999999, 0, 0, -1, 0.8)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F71:AC71").Value2 = $arr

$ws.Range("B72").Value2 = 6832552
$row = @("Sharjah SCC", "Khor Fakkan", 4, 1, "H", 1.3, 5.25, 7.5, 1.3, 5.25, 7.5, -1.75, 1.95, 1.85, 3.5, 1.975, 1.825, 0.3, -1, -1, 0.95, -1, 0.9750000000000001, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F72:AC72").Value2 = $arr

$ws.Range("B73").Value2 = 6832551
$row = @("Baniyas SC", "Shabab Al Ahli Dubai", 1, 2, "A", 5.25, 4.75, 1.45, 5, 4.75, 1.45, 1.25, 1.85, 1.95, 3.5, 1.8, 2, -1, -1, 0.45, 0.425, -0.5, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F73:AC73").Value2 = $arr

$ws.Range("B74").Value2 = 6832553
$row = @("Al Wasl SC", "Al Ittihad Kalba", 2, 1, "H", 1.5, 4.333, 5.25, 1.444, 4.5, 5.5, -1.25, 1.9, 1.9, 3.5, 1.975, 1.825, 0.444, -1, -1, -0.5, 0.45, -1, 0.825)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F74:AC74").Value2 = $arr

$ws.Range("B75").Value2 = 6832554
$row = @("Ajman SCC", "Al Nasr SC", 0, 3, "A", 3.4, 3.75, 1.85, 3.1, 3.6, 2, 0.25, 2.025, 1.775, 3, 2, 1.8, -1, -1, 1, -1, 0.7749999999999999, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F75:AC75").Value2 = $arr

$ws.Range("B79").Value2 = 6832559
$row = @("Al Nasr SC", "Emirates Club RAK", 2, 1, "H", 1.4, 4.75, 6, 1.333, 5, 7, -1.5, 1.9, 1.9, 3.5, 1.925, 1.875, 0.333, -1, -1, -1, 0.8999999999999999, -1, 0.875)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F79:AC79").Value2 = $arr

$ws.Range("B80").Value2 = 6832708
$row = @("Al Jazira SC", "Hatta Dubai", 3, 1, "H", 1.222, 6, 10, 1.222, 6, 9.5, -2, 1.9, 1.9, 4, 2, 1.8, 0.222, -1, -1, 0, 0, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F80:AC80").Value2 = $arr

$ws.Range("B91").Value2 = 6832566
$row = @("Al Jazira SC", "Al Nasr SC", 2, 3, "A", 1.533, 4.5, 4.333, 2.25, 4, 2.55, 0, 1.775, 2.025, 3.25, 2, 1.8, -1, -1, 1.55, -1, 1.025, 1, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F91:AC91").Value2 = $arr

$ws.Range("B92").Value2 = 6832567
$row = @("Sharjah SCC", "Al Ittihad Kalba", 1, 0, "H", 1.571, 4.333, 4.333, 1.6, 4.2, 4.5, -1, 2, 1.8, 3, 1.825, 1.975, 0.6000000000000001, -1, -1, 0, 0, -1, 0.9750000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt $row.Length; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("F92:AC92").Value2 = $arr

# --- New rows 96-102 (new fixtures added at bottom) ---
# Row 96
$ws.Range("A96").Value2 = 94
$ws.Range("B96").Value2 = 6832713
$ws.Range("C96").Value2 = "UAE Premier League"
$ws.Range("D96").Value2 = "UAE Premier League"
$ws.Range("E96").Value2 = (Get-Date -Year 2024 -Month 2 -Day 29 -Hour 10 -Minute 35 -Second 0)
$ws.Range("F96").Value2 = "Al Ain SCC"
$ws.Range("G96").Value2 = "Hatta Dubai"
$rowVals = @(1.125, 7.0, 15.0, 1.125, 7.5, 15.0, -2.5, 2.0, 1.8, 3.75, 1.975, 1.825, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K96:AA96").Value2 = $arr

# Row 97
$ws.Range("A97").Value2 = 95
$ws.Range("B97").Value2 = 6832571
$ws.Range("C97").Value2 = "UAE Premier League"
$ws.Range("D97").Value2 = "UAE Premier League"
$ws.Range("E97").Value2 = (Get-Date -Year 2024 -Month 2 -Day 29 -Hour 13 -Minute 15 -Second 0)
$ws.Range("F97").Value2 = "Al Ittihad Kalba"
$ws.Range("G97").Value2 = "Al Bataeh"
$rowVals = @(1.8, 3.8, 3.6, 1.8, 3.8, 3.6, -0.5, 1.8, 2.0, 3.0, 1.875, 1.925, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K97:AA97").Value2 = $arr

# Row 98
$ws.Range("A98").Value2 = 96
$ws.Range("B98").Value2 = 6832576
$ws.Range("C98").Value2 = "UAE Premier League"
$ws.Range("D98").Value2 = "UAE Premier League"
$ws.Range("E98").Value2 = (Get-Date -Year 2024 -Month 3 -Day 1 -Hour 10 -Minute 35 -Second 0)
$ws.Range("F98").Value2 = "Ajman SCC"
$ws.Range("G98").Value2 = "Emirates Club RAK"
$rowVals = @(1.6, 4.0, 4.75, 1.615, 4.0, 4.75, -0.75, 1.8, 2.0, 3.25, 1.85, 1.95, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K98:AA98").Value2 = $arr

# Row 99
$ws.Range("A99").Value2 = 97
$ws.Range("B99").Value2 = 6832572
$ws.Range("C99").Value2 = "UAE Premier League"
$ws.Range("D99").Value2 = "UAE Premier League"
$ws.Range("E99").Value2 = (Get-Date -Year 2024 -Month 3 -Day 1 -Hour 10 -Minute 35 -Second 0)
$ws.Range("F99").Value2 = "Khor Fakkan"
$ws.Range("G99").Value2 = "Shabab Al Ahli Dubai"
$rowVals = @(6.0, 5.5, 1.333, 5.25, 5.0, 1.444, 1.25, 1.9, 1.9, 3.5, 2.0, 1.8, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K99:AA99").Value2 = $arr

# Row 100
$ws.Range("A100").Value2 = 98
$ws.Range("B100").Value2 = 6832575
$ws.Range("C100").Value2 = "UAE Premier League"
$ws.Range("D100").Value2 = "UAE Premier League"
$ws.Range("E100").Value2 = (Get-Date -Year 2024 -Month 3 -Day 1 -Hour 13 -Minute 15 -Second 0)
$ws.Range("F100").Value2 = "Al Wasl SC"
$ws.Range("G100").Value2 = "Sharjah SCC"
$rowVals = @(1.85, 3.75, 3.5, 1.75, 4.0, 3.75, -0.75, 1.975, 1.825, 3.0, 1.8, 2.0, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K100:AA100").Value2 = $arr

# Row 101
$ws.Range("A101").Value2 = 99
$ws.Range("B101").Value2 = 6832574
$ws.Range("C101").Value2 = "UAE Premier League"
$ws.Range("D101").Value2 = "UAE Premier League"
$ws.Range("E101").Value2 = (Get-Date -Year 2024 -Month 3 -Day 2 -Hour 10 -Minute 35 -Second 0)
$ws.Range("F101").Value2 = "Al Nasr SC"
$ws.Range("G101").Value2 = "Baniyas SC"
$rowVals = @(1.571, 4.0, 4.75, 1.85, 3.6, 3.5, -0.5, 1.85, 1.95, 2.75, 1.8, 2.0, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K101:AA101").Value2 = $arr

# Row 102
$ws.Range("A102").Value2 = 100
$ws.Range("B102").Value2 = 6832573
$ws.Range("C102").Value2 = "UAE Premier League"
$ws.Range("D102").Value2 = "UAE Premier League"
$ws.Range("E102").Value2 = (Get-Date -Year 2024 -Month 3 -Day 2 -Hour 13 -Minute 15 -Second 0)
$ws.Range("F102").Value2 = "Al Wahda Abu Dhabi"
$ws.Range("G102").Value2 = "Al Jazira SC"
$rowVals = @(1.95, 3.75, 3.3, 2.05, 3.6, 3.2, -0.25, 1.775, 2.025, 3.25, 1.925, 1.875, 0.0, 0.0, 0.0, 0.0, 0.0)
$arr = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("K102:AA102").Value2 = $arr


# --- Formatting for the newly appended rows: copy row 95's look (bold/
#     bordered/centered column A style and the YYYY-MM-DD HH:MM:SS date
#     format used by column E) down onto rows 96-102. ---
$ws.Rows.Item(95).Copy() | Out-Null
$ws.Range("A96:AC102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
